$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.766.21"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "2.666.10"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.27"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.23"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.660"
$ws.Range("E7").Value = "  +6.96%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -3.58%  "
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.83"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("E13").Value = "  -2.19%  "
$ws.Range("E14").Value = "  -2.51%  "
$ws.Range("D15").Value = "3.140.95"
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").Value = "65.621.28"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").Value = "2.607.63"
$ws.Range("E17").Value = "  -2.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.57"
$ws.Range("E18").Value = "  -2.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.77"
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "349.83"
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.43"
$ws.Range("E21").Value = "  -3.29%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.66"
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("E24").Value = "  +8.83%  "
$ws.Range("E25").Value = "  +1.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.54"
$ws.Range("E26").Value = "  -2.13%  "
$ws.Range("E27").Value = "  +1.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "563.70"
$ws.Range("E28").Value = "  +5.88%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.163"
$ws.Range("E29").Value = "  -2.73%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.03"
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.14"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("E33").Value = "  +2.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.57"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("E35").Value = "  -0.67%  "
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.48"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "154.57"
$ws.Range("E40").Value = "  -2.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "160.71"
$ws.Range("E41").Value = "  -2.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.08"
$ws.Range("E42").Value = "  -0.96%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.28"
$ws.Range("E43").Value = "  -1.33%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0603"
$ws.Range("E44").Value = "  -1.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.68"
$ws.Range("E45").Value = "  -1.73%  "
$ws.Range("E46").Value = "  -0.72%  "
$ws.Range("E47").Value = "  +2.72%  "
$ws.Range("E48").Value = "  -1.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.78"
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("D50").Value = "0.0₆0244"
$ws.Range("E50").Value = "  +2.63%  "
$ws.Range("E51").Value = "  -1.90%  "
